$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in new task rows (5 thru 9) and extend existing rows (2 thru 4) ---

# Row 2 (Git Setup) - add Status/Completed On/Delay columns
$ws.Range("E2").Value = "Completed"
$ws.Range("F2").Value = "26/3/2017"
$ws.Range("G2").Value = "2 days"

# Row 3 (User story elaboration) - add Status
$ws.Range("E3").Value = "Yet to start"

# Row 4 (Database Designing) - add Status
$ws.Range("E4").Value = "Yet to start"

# Row 5 (Basic Blog)
$ws.Range("A5").Value = "Basic Blog"
$ws.Range("B5").Value = "Karan"
$ws.Range("C5").Value = "27/3/2017 "
$ws.Range("D5").Value = "27/3/2017 "
$ws.Range("E5").Value = "Completed"
$ws.Range("F5").Value = "27/3/2017 "
$ws.Range("G5").Value = "None"

# Row 6 (Basic Authentication)
$ws.Range("A6").Value = "Basic Authentication"
$ws.Range("B6").Value = "Karan"
$ws.Range("C6").Value = "27/3/2017 "
$ws.Range("D6").Value = "28/3/2017 "
$ws.Range("E6").Value = "Completed"
$ws.Range("F6").Value = "28/3/2017"
$ws.Range("G6").Value = "None"

# Row 7 (Basic Registration)
$ws.Range("A7").Value = "Basic Registration"
$ws.Range("B7").Value = "Karan"
$ws.Range("C7").Value = "27/3/2017"
$ws.Range("D7").Value = "28/3/2017"
$ws.Range("E7").Value = "In Progress"

# Row 8 (Basic searching)
$ws.Range("A8").Value = "Basic searching"
$ws.Range("B8").Value = "Karan"
$ws.Range("C8").Value = "27/3/2017 "
$ws.Range("D8").Value = "28/3/2017"

# Row 9 (Basic Forum)
$ws.Range("A9").Value = "Basic Forum"
$ws.Range("B9").Value = "Karan"
$ws.Range("C9").Value = "27/3/2017 "
$ws.Range("D9").Value = "30/3/2017"

# --- Formatting ---

# Thin border around the header row (merges with existing bold font style)
$ws.Range("A1:G1").Borders.LineStyle = 1

# Thin border around all the data rows that have content
$ws.Range("A2:D4").Borders.LineStyle = 1
$ws.Range("E2:E4").Borders.LineStyle = 1
$ws.Range("F2:F2").Borders.LineStyle = 1
$ws.Range("G2:G2").Borders.LineStyle = 1
$ws.Range("A5:D9").Borders.LineStyle = 1
$ws.Range("E5:E7").Borders.LineStyle = 1
$ws.Range("F5:F6").Borders.LineStyle = 1
$ws.Range("G5:G6").Borders.LineStyle = 1

# Highlight fills
$ws.Range("E2").Interior.Color = 65535
$ws.Range("E5").Interior.Color = 5296274
$ws.Range("E6").Interior.Color = 5296274

# --- Selection ---
$null = $ws.Range("E7").Select()
